# Recompute column G ("log Height sq") as the square of column F ("log Height").
# Previously G held 2*F (i.e. log(Height^2) simplified to 2*log(Height)); the
# corrected formula squares the log-height value instead: G = F^2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $logHeight = $ws.Cells.Item($r, 6).Value2
    if ($logHeight -ne $null) {
        $ws.Cells.Item($r, 7).Value = $logHeight * $logHeight
    }
}

# Restore the active selection to H9, matching the saved cursor position.
$ws.Range("H9").Select()
